# Add files via upload
# Update the "pie_threshold_range" Max value (C5) from 15 to 20, and move
# the active selection from B3 to C5 to match where the edit was made.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C5").Value = 20
$ws.Range("C5").Select()
